$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The second data row ("Manca_898_user") was actually stored as row 3 in the
# source file, leaving row 2 completely blank/unused. Deleting that blank
# row shifts every row below it up by one -- which is exactly the row
# layout the target file has (row 3 -> row 2, row 4 -> row 3, ...,
# row 13 -> row 12), with the same cell values/styles carried along.
$ws.Rows(2).Delete()

# Reflect the new cursor position left behind by the edit.
$ws.Range("A2").Select() | Out-Null
